$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-25 Wednesday", "2024-12-26 Thursday"),
    @("409÷5=", "676÷6="),
    @("137÷8=", "256÷3="),
    @("124÷8=", "602÷6="),
    @("785÷8=", "200÷3="),
    @("897÷9=", "462÷8="),
    @("489÷3=", "253÷9="),
    @("267÷2=", "452÷9="),
    @("116÷5=", "630÷4="),
    @("708÷6=", "520÷5="),
    @("296÷6=", "824÷5="),
    @("553÷2=", "139÷7="),
    @("940÷6=", "430÷8="),
    @("796÷8=", "333÷7="),
    @("690÷8=", "286÷8="),
    @("716÷8=", "949÷9="),
    @("252÷2=", "441÷4="),
    @("440÷9=", "859÷4="),
    @("409÷2=", "160÷7="),
    @("909÷6=", "648÷9="),
    @("822÷2=", "244÷4="),
    @("314÷5=", "500÷3="),
    @("783÷7=", "484÷4="),
    @("524÷2=", "656÷7="),
    @("485÷6=", "589÷2="),
    @("484÷9=", "479÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
